$d = $word.ActiveDocument

# 1) Main document body: the bold "TERE" in "A TERE," becomes "QWER".
$d.Content.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# 2) Header (default/primary header of section 1) contains several
#    "TRE"/"TERE"/"Tre"/"tre" placeholders that all get replaced too,
#    each with a (possibly different) "QWER"/"Qwer"/"Qewr"/"qwer" value.
#    They must be replaced in document order using one advancing Range
#    so that repeated identical substrings are matched positionally.
$hdr = $d.Sections.Item(1).Headers.Item(1)
$rng = $hdr.Range

$replacements = @(
    @("TRE", "QWER"),
    @("TERE", "QWER"),
    @("Tre", "Qwer"),
    @("Tre", "Qwer"),
    @("Tre", "Qewr"),
    @("Tre", "Qewr"),
    @("Tre", "Qwer"),
    @("tre", "qwer"),
    @("tre", "qwer"),
    @("tre", "qwer")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null
}
